# Added Week 15 simulations
$wb = $excel.ActiveWorkbook

# --- Sheet "OFF": row 3 (B3:G3) updated ---
$wsOff = $wb.Worksheets.Item("OFF")
$wsOff.Range("B3").Value = 325
$wsOff.Range("C3").Value = 214
$wsOff.Range("D3").Value = 58
$wsOff.Range("E3").Value = 21
$wsOff.Range("F3").Value = 7
$wsOff.Range("G3").Value = 7

# --- Sheet "DEF": row 3 (B3:E3) updated ---
$wsDef = $wb.Worksheets.Item("DEF")
$wsDef.Range("B3").Value = 414
$wsDef.Range("C3").Value = 298
$wsDef.Range("D3").Value = 87
$wsDef.Range("E3").Value = 45
